# ADD results from server
# Update row 2 data values on each year sheet with the latest results from the server run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("2025")
$ws.Cells.Item(2, "A").Value = 0
$ws.Cells.Item(2, "B").Value = 290.0628494009472
$ws.Cells.Item(2, "C").Value = 0
$ws.Cells.Item(2, "D").Value = 0
$ws.Cells.Item(2, "E").Value = 29049.07128553879
$ws.Cells.Item(2, "F").Value = 0
$ws.Cells.Item(2, "G").Value = 8095.92571266193
$ws.Cells.Item(2, "H").Value = 0
$ws.Cells.Item(2, "I").Value = 14940.21181152952
$ws.Cells.Item(2, "J").Value = 0
$ws.Cells.Item(2, "K").Value = 0
$ws.Cells.Item(2, "L").Value = 50998.86069102
$ws.Cells.Item(2, "M").Value = 11228.70813999
$ws.Cells.Item(2, "N").Value = 7234.066823234584
$ws.Cells.Item(2, "O").Value = 6709.085762003203

$ws = $wb.Worksheets.Item("2030")
$ws.Cells.Item(2, "A").Value = 219.6191807040655
$ws.Cells.Item(2, "B").Value = 3803.73674200606
$ws.Cells.Item(2, "C").Value = 0
$ws.Cells.Item(2, "D").Value = 0
$ws.Cells.Item(2, "E").Value = 45497.54863566629
$ws.Cells.Item(2, "F").Value = 0
$ws.Cells.Item(2, "G").Value = 8095.92571266193
$ws.Cells.Item(2, "H").Value = 0
$ws.Cells.Item(2, "I").Value = 31289.72026470282
$ws.Cells.Item(2, "J").Value = 0
$ws.Cells.Item(2, "K").Value = 0
$ws.Cells.Item(2, "L").Value = 60054.94214326091
$ws.Cells.Item(2, "M").Value = 17372.009741075
$ws.Cells.Item(2, "N").Value = 9195.862605783142
$ws.Cells.Item(2, "O").Value = 7881.258210869395

$ws = $wb.Worksheets.Item("2035")
$ws.Cells.Item(2, "A").Value = 2152.642195393625
$ws.Cells.Item(2, "B").Value = 5758.176568655313
$ws.Cells.Item(2, "C").Value = 0
$ws.Cells.Item(2, "D").Value = 0
$ws.Cells.Item(2, "E").Value = 57498.34502187894
$ws.Cells.Item(2, "F").Value = 0
$ws.Cells.Item(2, "G").Value = 8095.92571266193
$ws.Cells.Item(2, "H").Value = 0
$ws.Cells.Item(2, "I").Value = 48286.53531869316
$ws.Cells.Item(2, "J").Value = 0
$ws.Cells.Item(2, "K").Value = 0
$ws.Cells.Item(2, "L").Value = 60054.94214326091
$ws.Cells.Item(2, "M").Value = 23340.7500113645
$ws.Cells.Item(2, "N").Value = 13671.28338866544
$ws.Cells.Item(2, "O").Value = 13096.10328573595

$ws = $wb.Worksheets.Item("2040")
$ws.Cells.Item(2, "A").Value = 2152.642195393625
$ws.Cells.Item(2, "B").Value = 5758.176568655313
$ws.Cells.Item(2, "C").Value = 0
$ws.Cells.Item(2, "D").Value = 0
$ws.Cells.Item(2, "E").Value = 57498.34502187894
$ws.Cells.Item(2, "F").Value = 0
$ws.Cells.Item(2, "G").Value = 8095.92571266193
$ws.Cells.Item(2, "H").Value = 0
$ws.Cells.Item(2, "I").Value = 48286.53531869316
$ws.Cells.Item(2, "J").Value = 0
$ws.Cells.Item(2, "K").Value = 0
$ws.Cells.Item(2, "L").Value = 60054.94214326091
$ws.Cells.Item(2, "M").Value = 23340.7500113645
$ws.Cells.Item(2, "N").Value = 13671.28338866544
$ws.Cells.Item(2, "O").Value = 13096.10328573595

$ws = $wb.Worksheets.Item("2045")
$ws.Cells.Item(2, "A").Value = 2152.642195393625
$ws.Cells.Item(2, "B").Value = 5758.176568655313
$ws.Cells.Item(2, "C").Value = 0
$ws.Cells.Item(2, "D").Value = 0
$ws.Cells.Item(2, "E").Value = 57498.34502187894
$ws.Cells.Item(2, "F").Value = 0
$ws.Cells.Item(2, "G").Value = 8095.92571266193
$ws.Cells.Item(2, "H").Value = 0
$ws.Cells.Item(2, "I").Value = 48286.53531869316
$ws.Cells.Item(2, "J").Value = 0
$ws.Cells.Item(2, "K").Value = 0
$ws.Cells.Item(2, "L").Value = 60054.94214326091
$ws.Cells.Item(2, "M").Value = 23340.7500113645
$ws.Cells.Item(2, "N").Value = 13671.28338866544
$ws.Cells.Item(2, "O").Value = 13096.10328573595

$ws = $wb.Worksheets.Item("2050")
$ws.Cells.Item(2, "A").Value = 2152.642195393625
$ws.Cells.Item(2, "B").Value = 5758.176568655313
$ws.Cells.Item(2, "C").Value = 0
$ws.Cells.Item(2, "D").Value = 0
$ws.Cells.Item(2, "E").Value = 57498.34502187894
$ws.Cells.Item(2, "F").Value = 0
$ws.Cells.Item(2, "G").Value = 8095.92571266193
$ws.Cells.Item(2, "H").Value = 0
$ws.Cells.Item(2, "I").Value = 48286.53531869316
$ws.Cells.Item(2, "J").Value = 0
$ws.Cells.Item(2, "K").Value = 0
$ws.Cells.Item(2, "L").Value = 60054.94214326091
$ws.Cells.Item(2, "M").Value = 23340.7500113645
$ws.Cells.Item(2, "N").Value = 13671.28338866544
$ws.Cells.Item(2, "O").Value = 13096.10328573595

